# Update existing row 37's timestamp (sub-millisecond correction) and
# append rows 38-43 with new price observations, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37: only the date/time value changes (tiny precision correction).
$ws.Range("A37").Value = 45833.45891689815

# New rows 38-43: date (col A), product (col B), weight (col C), price (col D).
$newRows = @(
    @{ Row = 38; Fecha = 45834.45434878472; Producto = "CREATINA MONOHIDRATO EN POLVO"; Peso = "1Kg"; Precio = "12,88€" },
    @{ Row = 39; Fecha = 45853.38054951389; Producto = "CREATINA MONOHIDRATO EN POLVO"; Peso = "1Kg"; Precio = "15,41€" },
    @{ Row = 40; Fecha = 45853.38514180556; Producto = "CREATINA MONOHIDRATO EN POLVO"; Peso = "1Kg"; Precio = "15,41€" },
    @{ Row = 41; Fecha = 45853.38685481482; Producto = "CREATINA MONOHIDRATO EN POLVO"; Peso = "1Kg"; Precio = "15,41€" },
    @{ Row = 42; Fecha = 45853.39298033565; Producto = "CREATINA MONOHIDRATO EN POLVO"; Peso = "1Kg"; Precio = "15,41€" },
    @{ Row = 43; Fecha = 45853.39427826957; Producto = "CREATINA MONOHIDRATO EN POLVO"; Peso = "1Kg"; Precio = "15,41€" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Column A keeps the same date/time number format used by the other
    # rows in this column (style index 2 in the original workbook).
    $ws.Range("A$rowNum").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("A$rowNum").Value = $r.Fecha

    $ws.Range("B$rowNum").Value = $r.Producto
    $ws.Range("C$rowNum").Value = $r.Peso
    $ws.Range("D$rowNum").Value = $r.Precio
}
